$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (357) down to the new rows (358-366)
# so the new date cells in column A keep the same number format / font / borders / alignment.
$ws.Range("A357:D357").Copy()
$ws.Range("A358:D366").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New daily data, continuing the series through 1 September 2021 (serial 44440)
$newData = @(
    @(44432, 1, 7, 583.8198498748957),
    @(44433, 0, 7, 583.8198498748957),
    @(44434, 0, 5, 417.0141784820684),
    @(44435, 4, 7, 583.8198498748957),
    @(44436, 2, 9, 750.6255212677231),
    @(44437, 1, 9, 750.6255212677231),
    @(44438, 0, 8, 667.2226855713094),
    @(44439, 0, 7, 583.8198498748957),
    @(44440, 0, 7, 583.8198498748957)
)

$row = 358
foreach ($rec in $newData) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row = $row + 1
}
